# Add a small bird-sighting table to the sheet: a bold, bordered,
# centered/top-aligned header row followed by two plain data rows.
# (commit message: "this lets it see the curnt date")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headers = @("Bird", "Number", "Location", "When")
for ($c = 1; $c -le $headers.Length; $c++) {
    $ws.Cells.Item(1, $c).Value = $headers[$c - 1]
}

# Build the header style once on A1, then propagate it to B1:D1 with a
# single PasteSpecial (formats-only) so every header cell ends up sharing
# exactly one cell style, instead of each individual property assignment
# minting its own transient style entry.
$a1 = $ws.Cells.Item(1, 1)
$a1.Font.Bold = $true
$a1.HorizontalAlignment = -4108   # xlCenter
$a1.VerticalAlignment = -4160     # xlTop
$a1.Borders.LineStyle = 1         # xlContinuous (thin)

$a1.Copy()
$ws.Range("B1:D1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$data = @(
    @("kiwi", "2", "50 hicks rd", "20/8/2025"),
    @("kiwi", "1", "airstrip", "21/8/2025")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c).Value = $row[$c - 1]
    }
}
